$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column BJ (14-aug) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy header style (bold, border, centered) from BI1 into the new BJ1 cell
$ws1.Cells.Item(1, 61).Copy($ws1.Cells.Item(1, 62))
$ws1.Cells.Item(1, 62).Value = "14-aug"

$ws1.Cells.Item(2, 62).Value = 98.09999999999999
$ws1.Cells.Item(3, 62).Value = 92
$ws1.Cells.Item(4, 62).Value = 88.45
$ws1.Cells.Item(5, 62).Value = 81.83
$ws1.Cells.Item(6, 62).Value = 82.70999999999999
$ws1.Cells.Item(7, 62).Value = 79.59
$ws1.Cells.Item(8, 62).Value = 84.89
$ws1.Cells.Item(9, 62).Value = 103.11
$ws1.Cells.Item(10, 62).Value = 99.31
$ws1.Cells.Item(11, 62).Value = 90.92
$ws1.Cells.Item(12, 62).Value = 76.66
$ws1.Cells.Item(13, 62).Value = 67.40000000000001
$ws1.Cells.Item(14, 62).Value = 54.01
$ws1.Cells.Item(15, 62).Value = 40.91
$ws1.Cells.Item(16, 62).Value = 54.9
$ws1.Cells.Item(17, 62).Value = 69.43000000000001
$ws1.Cells.Item(18, 62).Value = 78.98999999999999
$ws1.Cells.Item(19, 62).Value = 84.63
$ws1.Cells.Item(20, 62).Value = 102
$ws1.Cells.Item(21, 62).Value = 115.78
$ws1.Cells.Item(22, 62).Value = 137.99
$ws1.Cells.Item(23, 62).Value = 142.01
$ws1.Cells.Item(24, 62).Value = 120.75
$ws1.Cells.Item(25, 62).Value = 105.95

# --- Sheet "Gaz": append row 59 (2025-08-12) ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Cells.Item(59, 1).NumberFormat = "@"
$ws2.Cells.Item(59, 1).Value = "2025-08-12"
$ws2.Cells.Item(59, 2).Value = 31.225

# --- Sheet "CO2": append row 59 (2025-08-12) ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Cells.Item(59, 1).NumberFormat = "@"
$ws3.Cells.Item(59, 1).Value = "2025-08-12"
$ws3.Cells.Item(59, 2).Value = 70.85
